$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# "Generate Report for Handback": the 57ca0734-... file's handback
# transform failed, so its status flips from "Ready for handoff" to
# "Handback transform failed", and the per-locale sheets now carry an
# Error Detail message explaining the mismatched file name.
# -----------------------------------------------------------------

$newStatus = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("P3").Value = "Handback file name: xc1dni2s.jyd is different with handoff file name: 57ca0734-c4b0-4dae-b279-8ec6d4ab8b2a.61c3f1a5f73f621f2a86a7fe20a231c7e90aa040.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("P3").Value = "Handback file name: xc1dni2s.jyd is different with handoff file name: 57ca0734-c4b0-4dae-b279-8ec6d4ab8b2a.61c3f1a5f73f621f2a86a7fe20a231c7e90aa040.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
